# Rename the "FP" worksheet to "Comprar um carro" (the comparison-matrix
# sheet now shows the criteria-comparison matrix's name on its tab).
$wb = $excel.ActiveWorkbook

$ws = $null
foreach ($sheet in $wb.Worksheets) {
    if ($sheet.Name -eq "FP") {
        $ws = $sheet
        break
    }
}

if ($ws -eq $null) {
    $ws = $wb.Worksheets.Item(6)
}

$ws.Name = "Comprar um carro"
